$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83, shifting existing rows 83..165 down to 84..166.
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new record's data.
$ws.Cells.Item(83, 1).Value = 6
$ws.Cells.Item(83, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(83, 3).Value = "Metropolitana"
$ws.Cells.Item(83, 4).Value = 44482
$ws.Cells.Item(83, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(83, 5).Value = 13
$ws.Cells.Item(83, 6).Value = "Fruta"
$ws.Cells.Item(83, 7).Value = 100101
$ws.Cells.Item(83, 8).Value = "Berries"
$ws.Cells.Item(83, 9).Value = 100101001
$ws.Cells.Item(83, 10).Value = "Arándano (blue)"
$ws.Cells.Item(83, 11).Value = "Sin especificar"
$ws.Cells.Item(83, 12).Value = "Especial"
$ws.Cells.Item(83, 13).Value = 750
$ws.Cells.Item(83, 14).Value = 13000
$ws.Cells.Item(83, 15).Value = 13000
$ws.Cells.Item(83, 16).Value = 13000
$ws.Cells.Item(83, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(83, 18).Value = "Perú"
$ws.Cells.Item(83, 19).Value = 6500
$ws.Cells.Item(83, 20).Value = 2
